# Fill in the student's name and personal identifier on the checklist,
# mark two more "Must Haves"/bonus items as satisfied (adds a score of
# 1 and 0.5 respectively, which the existing Sum Points formula in B54
# picks up automatically), and leave the active selection on F7 (instead
# of the previous scrolled-down E42 view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 / A5 held bare "Student Name:" / "Personal Identifier:" labels —
# fill in the actual values.
$ws.Range("A4").Value = "Student Name: Brian Schneider"
$ws.Range("A5").Value = "Personal Identifier: if21b072"

# Row 33 ("consider specialities in battle-rounds between cards") and
# row 52 ("Contains link to GIT") now score points in column B.
$ws.Range("B33").Value = 1
$ws.Range("B52").Value = 0.5

# Restore the view/selection to F7 (B54's SUM formula recalculates on
# its own once B33/B52 are populated).
$ws.Range("F7").Select() | Out-Null
